$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell without letting Excel's
# autodetection reinterpret a date-looking string ("2023-09-07") as a
# real date serial. Flip the format to text, assign, then restore the
# original style/format so nothing else about the cell changes.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# --- Row 10 / Row 11: swap the survey-point identity fields (A, P, Q, R) ---
$a10 = $ws.Cells.Item(10, 1).Value2
$p10 = $ws.Cells.Item(10, 16).Value2
$q10 = $ws.Cells.Item(10, 17).Value2
$r10 = $ws.Cells.Item(10, 18).Value2

$a11 = $ws.Cells.Item(11, 1).Value2
$p11 = $ws.Cells.Item(11, 16).Value2
$q11 = $ws.Cells.Item(11, 17).Value2
$r11 = $ws.Cells.Item(11, 18).Value2

$ws.Cells.Item(10, 1).Value = $a11
$ws.Cells.Item(10, 16).Value = $p11
$ws.Cells.Item(10, 17).Value = $q11
$ws.Cells.Item(10, 18).Value = $r11

$ws.Cells.Item(11, 1).Value = $a10
$ws.Cells.Item(11, 16).Value = $p10
$ws.Cells.Item(11, 17).Value = $q10
$ws.Cells.Item(11, 18).Value = $r10

# --- Row 14: replace the record with the Kolflarnlav observation ---
$ws.Cells.Item(14, 1).Value = 111942712
$ws.Cells.Item(14, 2).Value = 77267
$ws.Cells.Item(14, 4).Value = "NT"
$ws.Cells.Item(14, 5).Value = 6446
$ws.Cells.Item(14, 6).Value = "Kolflarnlav"
$ws.Cells.Item(14, 7).Value = "Carbonicola anthracophila"
$ws.Cells.Item(14, 8).Value = "(Nyl.) Bendiksby & Timdal"
$ws.Cells.Item(14, 9).Value = ""
$ws.Cells.Item(14, 17).Value = 468231.4750461024
$ws.Cells.Item(14, 18).Value = 6875021.661872049

Set-TextValue $ws.Cells.Item(14, 25) "2023-09-07"
Set-TextValue $ws.Cells.Item(14, 27) "2023-09-07"

# --- Drop the now-obsolete rows 15-19 (their data moved into row 14 above) ---
$ws.Range("A15:A19").EntireRow.Delete()
